$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title strings (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# --- Fix cells that switch from the literal "0"/"***.*" placeholder text to real numbers ---
# (copy number style from a sibling cell, then overwrite with the numeric value)
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 100
$ws.Range("C15").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("H15").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 100
$ws.Range("C15").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = 1
$ws.Range("H15").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = -100
$ws.Range("C15").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("G33").Value = 1
$ws.Range("H15").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = -100
$ws.Range("C15").Copy()
$ws.Range("J33").PasteSpecial(-4122)
$ws.Range("J33").Value = 1
$ws.Range("H15").Copy()
$ws.Range("K33").PasteSpecial(-4122)
$ws.Range("K33").Value = -100

# --- Fix cells that switch from a real number to the literal "0" placeholder text ---
# (copy value+format from a cell that already holds that exact placeholder)
$ws.Range("C23").Copy()
$ws.Range("C14").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C28").PasteSpecial(-4104)
$ws.Range("C23").Copy()
$ws.Range("C28").PasteSpecial(-4122)

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 2
$ws.Range("N14").Value = 100
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 150
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 20
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = -61.290322580645
$ws.Range("I16").Value = 15
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = -55.882352941176
$ws.Range("L16").Value = -34.782608695652
$ws.Range("M16").Value = -42.307692307692
$ws.Range("N16").Value = -88
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 43
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 2.380952380952
$ws.Range("L17").Value = -10.416666666666
$ws.Range("M17").Value = 104.761904761905
$ws.Range("N17").Value = 22.857142857142
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = -54.545454545454
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = -58.333333333333
$ws.Range("L18").Value = -16.666666666666
$ws.Range("M18").Value = -62.962962962963
$ws.Range("N18").Value = -96.254681647940
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -46.153846153846
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 99
$ws.Range("H19").Value = -48.484848484848
$ws.Range("I19").Value = 55
$ws.Range("J19").Value = 110
$ws.Range("K19").Value = -50
$ws.Range("L19").Value = -6.779661016949
$ws.Range("M19").Value = 14.583333333333
$ws.Range("N19").Value = -49.541284403669
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -38.095238095238
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = -40
$ws.Range("L20").Value = -34.782608695652
$ws.Range("M20").Value = -31.818181818181
$ws.Range("N20").Value = -92.385786802030
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 59
$ws.Range("E21").Value = -44.067796610169
$ws.Range("F21").Value = 130
$ws.Range("G21").Value = 212
$ws.Range("H21").Value = -38.679245283018
$ws.Range("I21").Value = 146
$ws.Range("J21").Value = 238
$ws.Range("K21").Value = -38.655462184873
$ws.Range("L21").Value = -13.095238095238
$ws.Range("M21").Value = -0.680272108843
$ws.Range("N21").Value = -80.243572395128
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -55.555555555555
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -44.444444444444
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 54
$ws.Range("E24").Value = -44.444444444444
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 197
$ws.Range("H24").Value = -44.162436548223
$ws.Range("I24").Value = 122
$ws.Range("J24").Value = 220
$ws.Range("K24").Value = -44.545454545454
$ws.Range("L24").Value = -39
$ws.Range("M24").Value = 27.083333333333
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = -29.629629629629
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 117
$ws.Range("H25").Value = -53.846153846153
$ws.Range("I25").Value = 58
$ws.Range("J25").Value = 129
$ws.Range("K25").Value = -55.038759689922
$ws.Range("L25").Value = -53.225806451612
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = 9.523809523809
$ws.Range("F26").Value = 85
$ws.Range("G26").Value = 64
$ws.Range("H26").Value = 32.8125
$ws.Range("I26").Value = 102
$ws.Range("J26").Value = 83
$ws.Range("K26").Value = 22.891566265060
$ws.Range("L26").Value = 45.714285714285
$ws.Range("M26").Value = 29.113924050632
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -12.5
$ws.Range("L27").Value = 40
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -60
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = -45.454545454545
$ws.Range("L28").Value = -25
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 2

$excel.CutCopyMode = 0
